$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto values. Values are prefixed with a literal
# apostrophe to force Excel to store them as text (matching the
# original inline-string cell type) instead of auto-converting
# numeric-looking strings (e.g. "1.003") into floating point numbers.
# The style is reset to "Normal" afterwards so no stray text
# number-format style gets attached to the cell.
function Set-TextValue($range, [string]$value) {
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "28.607.62"
Set-TextValue $ws.Range("D3") "1.791.77"
Set-TextValue $ws.Range("E3") "  -2.06%  "
Set-TextValue $ws.Range("D4") "1.003"
Set-TextValue $ws.Range("E4") "  +0.09%  "
Set-TextValue $ws.Range("D5") "231.62"
Set-TextValue $ws.Range("E5") "  -1.55%  "
Set-TextValue $ws.Range("D6") "0.5884"
Set-TextValue $ws.Range("E6") "  -2.40%  "
Set-TextValue $ws.Range("E7") "  +0.08%  "
Set-TextValue $ws.Range("D8") "0.2763"
Set-TextValue $ws.Range("E8") "  -1.06%  "
Set-TextValue $ws.Range("B9") "Dogecoin"
Set-TextValue $ws.Range("C9") "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextValue $ws.Range("D9") "0.06738"
Set-TextValue $ws.Range("E9") "  -4.53%  "
Set-TextValue $ws.Range("B10") "Solana"
Set-TextValue $ws.Range("C10") "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue $ws.Range("D10") "23.14"
Set-TextValue $ws.Range("E10") "  -1.91%  "
Set-TextValue $ws.Range("D11") "0.07531"
Set-TextValue $ws.Range("E11") "  -1.67%  "
Set-TextValue $ws.Range("D12") "1.793.24"
Set-TextValue $ws.Range("E12") "  -2.30%  "
Set-TextValue $ws.Range("D14") "0.6133"
Set-TextValue $ws.Range("E14") "  -2.42%  "
Set-TextValue $ws.Range("D15") "2.034.46"
Set-TextValue $ws.Range("E15") "  -2.12%  "
Set-TextValue $ws.Range("D16") "75.32"
Set-TextValue $ws.Range("E16") "  -4.76%  "
Set-TextValue $ws.Range("D17") "0.000008904"
Set-TextValue $ws.Range("E17") "  -9.84%  "
Set-TextValue $ws.Range("D18") "28.581.48"
Set-TextValue $ws.Range("E18") "  -2.22%  "
Set-TextValue $ws.Range("D19") "5.434"
Set-TextValue $ws.Range("E19") "  -6.93%  "
Set-TextValue $ws.Range("D21") "209.74"
Set-TextValue $ws.Range("E21") "  -6.47%  "
Set-TextValue $ws.Range("D22") "11.46"
Set-TextValue $ws.Range("E22") "  -2.10%  "
Set-TextValue $ws.Range("D23") "6.818"
Set-TextValue $ws.Range("E23") "  -2.77%  "
Set-TextValue $ws.Range("D24") "1.004"
Set-TextValue $ws.Range("E24") "  +0.10%  "
Set-TextValue $ws.Range("D25") "152.81"
Set-TextValue $ws.Range("E25") "  -2.31%  "
Set-TextValue $ws.Range("D26") "8.088"
Set-TextValue $ws.Range("E26") "  +1.38%  "
Set-TextValue $ws.Range("D27") "0.1259"
Set-TextValue $ws.Range("E27") "  -3.46%  "
Set-TextValue $ws.Range("D28") "16.38"
Set-TextValue $ws.Range("E28") "  -1.40%  "
Set-TextValue $ws.Range("D29") "1.417"
Set-TextValue $ws.Range("E29") "  -4.15%  "
Set-TextValue $ws.Range("D30") "0.06178"
Set-TextValue $ws.Range("E30") "  -5.08%  "
Set-TextValue $ws.Range("E31") "  -1.93%  "
Set-TextValue $ws.Range("D32") "3.808"
Set-TextValue $ws.Range("E32") "  +0.32%  "
Set-TextValue $ws.Range("D33") "3.784"
Set-TextValue $ws.Range("E33") "  -1.42%  "
Set-TextValue $ws.Range("D34") "1.736"
Set-TextValue $ws.Range("E34") "  +0.68%  "
Set-TextValue $ws.Range("D35") "1.046"
Set-TextValue $ws.Range("E35") "  -5.56%  "
Set-TextValue $ws.Range("D36") "0.6389"
Set-TextValue $ws.Range("E36") "  -1.29%  "
Set-TextValue $ws.Range("D37") "2.501"
Set-TextValue $ws.Range("E37") "  -1.71%  "
Set-TextValue $ws.Range("E38") "  -0.97%  "
Set-TextValue $ws.Range("D39") "6.408"
Set-TextValue $ws.Range("E39") "  -2.59%  "
Set-TextValue $ws.Range("D40") "0.01693"
Set-TextValue $ws.Range("E40") "  -3.10%  "
Set-TextValue $ws.Range("D41") "1.140.22"
Set-TextValue $ws.Range("E41") "  -6.16%  "
Set-TextValue $ws.Range("D42") "0.8769"
Set-TextValue $ws.Range("E42") "  -2.20%  "
Set-TextValue $ws.Range("E43") "  +0.29%  "
Set-TextValue $ws.Range("D44") "99.98"
Set-TextValue $ws.Range("E44") "  -0.40%  "
Set-TextValue $ws.Range("D45") "1.943.48"
Set-TextValue $ws.Range("E45") "  -2.65%  "
Set-TextValue $ws.Range("D46") "59.88"
Set-TextValue $ws.Range("E46") "  -4.47%  "
Set-TextValue $ws.Range("D47") "0.00000000110"
Set-TextValue $ws.Range("E47") "  -4.59%  "
Set-TextValue $ws.Range("D48") "1.585"
Set-TextValue $ws.Range("E48") "  +0.11%  "
Set-TextValue $ws.Range("D49") "8.363"
Set-TextValue $ws.Range("E49") "  -2.38%  "
Set-TextValue $ws.Range("D50") "0.05478"
Set-TextValue $ws.Range("E50") "  -0.39%  "
Set-TextValue $ws.Range("D51") "0.4475"
Set-TextValue $ws.Range("E51") "  -1.61%  "
